# Apply the edit described by the diff:
# - Add two new shared strings: "hintTooltip" and "Press this button to show hints."
# - Add a new row 22 to sheet "en": A22 = "hintTooltip", B22 = "Press this button to show hints." (wrap-text style)
# - Update the sheet selection to B22 (handled automatically by setting values / selecting range)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# New key/value pair appended at the bottom of the table (row 22)
$ws.Range("A22").Value = "hintTooltip"
$ws.Range("B22").Value = "Press this button to show hints."

# Match the wrap-text style used by the other "value" cells in column B (e.g. B21)
$ws.Range("B22").WrapText = $true

# Move the active selection to the newly added cell, matching the diff's sheetView selection
$ws.Range("B22").Select()

$wb.Save()
